# Auto-generated Excel COM-interop script to apply odds updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.56
$ws.Range("G2").Value = 1.65
$ws.Range("H2").Value = 7.6
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 3.75
$ws.Range("K2").Value = 4.2
$ws.Range("N2").Value = 2.88
$ws.Range("O2").Value = 1.44
$ws.Range("P2").Value = 1.62
$ws.Range("Q2").Value = 2.34
$ws.Range("R2").Value = 1.22
$ws.Range("T2").Value = 2.28
$ws.Range("U2").Value = 1.66
$ws.Range("V2").Value = 1.12
$ws.Range("W2").Value = 2.52
$ws.Range("X2").Value = 11.5
$ws.Range("Y2").Value = 1000
$ws.Range("AC2").Value = 9.4
$ws.Range("AE2").Value = 220
$ws.Range("AI2").Value = 210
$ws.Range("AJ2").Value = 18
$ws.Range("AM2").Value = 320

# Row 3
$ws.Range("L3").Value = 1.23
$ws.Range("S3").Value = 2.32
$ws.Range("U3").Value = 2.58
$ws.Range("AK3").Value = 48

# Row 4
$ws.Range("H4").Value = 1.9
$ws.Range("I4").Value = 1.91
$ws.Range("Q4").Value = 1.81
$ws.Range("S4").Value = 3.05
$ws.Range("U4").Value = 2.26
$ws.Range("V4").Value = 2.08
$ws.Range("AB4").Value = 18.5
$ws.Range("AH4").Value = 17.5

# Row 5
$ws.Range("Q5").Value = 1.66
$ws.Range("S5").Value = 2.62
$ws.Range("T5").Value = 1.59
$ws.Range("U5").Value = 2.62
$ws.Range("X5").Value = 20
$ws.Range("AB5").Value = 18

# Row 6
$ws.Range("F6").Value = 2.24
$ws.Range("G6").Value = 2.4
$ws.Range("L6").Value = 1.32
$ws.Range("T6").Value = 1.69
$ws.Range("W6").Value = 1.73

# Row 7
$ws.Range("F7").Value = 1.6
$ws.Range("H7").Value = 2.38
$ws.Range("J7").Value = 1.58
$ws.Range("V7").Value = 1.56

# Row 8
$ws.Range("G8").Value = 3.1
$ws.Range("J8").Value = 3.2
$ws.Range("L8").Value = 1.32
$ws.Range("T8").Value = 1.66
$ws.Range("U8").Value = 2.26

# Row 9
$ws.Range("K9").Value = 3.05
$ws.Range("L9").Value = 1.66
$ws.Range("Q9").Value = 2.96

# Row 11
$ws.Range("F11").Value = 3.25
$ws.Range("H11").Value = 2.56
$ws.Range("J11").Value = 2.9
$ws.Range("K11").Value = 3

# Row 12
$ws.Range("J12").Value = 2.88
$ws.Range("K12").Value = 3.15
$ws.Range("L12").Value = 1.65
$ws.Range("N12").Value = 2.32
$ws.Range("W12").Value = 1.84

# Row 13
$ws.Range("L13").Value = 1.77
$ws.Range("T13").Value = 2.5

# Row 15
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2.12
$ws.Range("I15").Value = 5.1
$ws.Range("J15").Value = 3.05
$ws.Range("K15").Value = 3.35
$ws.Range("L15").Value = 1.44
$ws.Range("M15").Value = 1.11
$ws.Range("N15").Value = 2.68
$ws.Range("O15").Value = 1.5
$ws.Range("P15").Value = 1.56
$ws.Range("Q15").Value = 2.56
$ws.Range("R15").Value = 1.21
$ws.Range("S15").Value = 5
$ws.Range("T15").Value = 2.08
$ws.Range("U15").Value = 1.75
$ws.Range("X15").Value = 11.5
$ws.Range("Y15").Value = 15
$ws.Range("Z15").Value = 970
$ws.Range("AB15").Value = 7.6
$ws.Range("AC15").Value = 9
$ws.Range("AF15").Value = 12
$ws.Range("AG15").Value = 11.5
$ws.Range("AH15").Value = 26
$ws.Range("AJ15").Value = 30
$ws.Range("AL15").Value = 55

# Row 16
$ws.Range("H16").Value = 2.44
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1.59

# Row 17
$ws.Range("F17").Value = 2.38
$ws.Range("G17").Value = 2.66
$ws.Range("I17").Value = 4.2
$ws.Range("J17").Value = 2.86
$ws.Range("K17").Value = 3.05
$ws.Range("M17").Value = 1.14
$ws.Range("N17").Value = 2.4
$ws.Range("O17").Value = 1.61
$ws.Range("P17").Value = 1.45
$ws.Range("Q17").Value = 2.64
$ws.Range("R17").Value = 1.16
$ws.Range("S17").Value = 5.5
$ws.Range("T17").Value = 2.16
$ws.Range("U17").Value = 1.68
$ws.Range("V17").Value = 1.31
$ws.Range("W17").Value = 1.6
$ws.Range("X17").Value = 8
$ws.Range("Y17").Value = 10.5
$ws.Range("Z17").Value = 970
$ws.Range("AB17").Value = 7.4
$ws.Range("AC17").Value = 7.4
$ws.Range("AD17").Value = 18.5
$ws.Range("AE17").Value = 75
$ws.Range("AF17").Value = 15
$ws.Range("AG17").Value = 13.5
$ws.Range("AH17").Value = 27
$ws.Range("AI17").Value = 120
$ws.Range("AJ17").Value = 970
$ws.Range("AK17").Value = 970
$ws.Range("AL17").Value = 75
$ws.Range("AM17").Value = 260
$ws.Range("AN17").Value = 55

# Row 18
$ws.Range("H18").Value = 4.3
$ws.Range("K18").Value = 3.2
$ws.Range("P18").Value = 1.58
$ws.Range("R18").Value = 1.21
$ws.Range("U18").Value = 1.78
